$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = -17.52614769222031
$ws.Cells.Item(2, 3).Value = -0.4794819392043819
$ws.Cells.Item(2, 4).Value = -17.52614769222031
$ws.Cells.Item(2, 5).Value = -17.52614769222031
$ws.Cells.Item(2, 6).Value = -17.52614769222031
$ws.Cells.Item(2, 7).Value = -17.52614769222031
$ws.Cells.Item(2, 8).Value = -17.52614769222031
$ws.Cells.Item(2, 9).Value = -17.52614769222031
$ws.Cells.Item(2, 10).Value = -17.52614769222031
$ws.Cells.Item(2, 11).Value = -17.52614769222031

# Row 3
$ws.Cells.Item(3, 2).Value = -17.52614769222031
$ws.Cells.Item(3, 3).Value = -17.52614769222031
$ws.Cells.Item(3, 4).Value = -17.52614769222031
$ws.Cells.Item(3, 5).Value = -17.52614769222031
$ws.Cells.Item(3, 6).Value = -17.52614769222031
$ws.Cells.Item(3, 7).Value = -17.52614769222031
$ws.Cells.Item(3, 8).Value = -17.52614769222031
$ws.Cells.Item(3, 9).Value = 0.620545272994536
$ws.Cells.Item(3, 10).Value = -17.52614769222031
$ws.Cells.Item(3, 11).Value = -17.52614769222031

# Row 4
$ws.Cells.Item(4, 2).Value = -17.52614769222031
$ws.Cells.Item(4, 3).Value = -0.3129959962955359
$ws.Cells.Item(4, 4).Value = 0.106968159785989
$ws.Cells.Item(4, 5).Value = -17.52614769222031
$ws.Cells.Item(4, 6).Value = 3.941404625821051
$ws.Cells.Item(4, 7).Value = -17.52614769222031
$ws.Cells.Item(4, 8).Value = 1.430821108887289
$ws.Cells.Item(4, 9).Value = -17.52614769222031
$ws.Cells.Item(4, 10).Value = 2.670450832555175
$ws.Cells.Item(4, 11).Value = -17.52614769222031

# Row 5
$ws.Cells.Item(5, 2).Value = -17.52614769222031
$ws.Cells.Item(5, 3).Value = 0.4596775715283858
$ws.Cells.Item(5, 4).Value = -17.52614769222031
$ws.Cells.Item(5, 5).Value = -17.52614769222031
$ws.Cells.Item(5, 6).Value = -17.52614769222031
$ws.Cells.Item(5, 7).Value = 3.406272537097631
$ws.Cells.Item(5, 8).Value = -17.52614769222031
$ws.Cells.Item(5, 9).Value = -17.52614769222031
$ws.Cells.Item(5, 10).Value = -17.52614769222031
$ws.Cells.Item(5, 11).Value = -17.52614769222031

# Row 6
$ws.Cells.Item(6, 2).Value = -17.52614769222031
$ws.Cells.Item(6, 3).Value = -17.52614769222031
$ws.Cells.Item(6, 4).Value = -17.52614769222031
$ws.Cells.Item(6, 5).Value = -17.52614769222031
$ws.Cells.Item(6, 6).Value = -17.52614769222031
$ws.Cells.Item(6, 7).Value = -17.52614769222031
$ws.Cells.Item(6, 8).Value = -17.52614769222031
$ws.Cells.Item(6, 9).Value = -17.52614769222031
$ws.Cells.Item(6, 10).Value = -17.52614769222031
$ws.Cells.Item(6, 11).Value = -17.52614769222031

# Row 7
$ws.Cells.Item(7, 2).Value = 3.237201665234361
$ws.Cells.Item(7, 3).Value = -17.52614769222031
$ws.Cells.Item(7, 4).Value = -17.52614769222031
$ws.Cells.Item(7, 5).Value = -17.52614769222031
$ws.Cells.Item(7, 6).Value = -17.52614769222031
$ws.Cells.Item(7, 7).Value = -17.52614769222031
$ws.Cells.Item(7, 8).Value = -17.52614769222031
$ws.Cells.Item(7, 9).Value = -17.52614769222031
$ws.Cells.Item(7, 10).Value = -17.52614769222031
$ws.Cells.Item(7, 11).Value = -17.52614769222031

# Row 8
$ws.Cells.Item(8, 2).Value = -17.52614769222031
$ws.Cells.Item(8, 3).Value = -17.52614769222031
$ws.Cells.Item(8, 4).Value = -17.52614769222031
$ws.Cells.Item(8, 5).Value = 1.231453918778268
$ws.Cells.Item(8, 6).Value = -17.52614769222031
$ws.Cells.Item(8, 7).Value = -17.52614769222031
$ws.Cells.Item(8, 8).Value = -17.52614769222031
$ws.Cells.Item(8, 9).Value = -17.52614769222031
$ws.Cells.Item(8, 10).Value = -17.52614769222031
$ws.Cells.Item(8, 11).Value = -17.52614769222031

# Row 9
$ws.Cells.Item(9, 2).Value = 3.401940508743195
$ws.Cells.Item(9, 3).Value = -17.52614769222031
$ws.Cells.Item(9, 4).Value = -17.52614769222031
$ws.Cells.Item(9, 5).Value = -17.52614769222031
$ws.Cells.Item(9, 6).Value = -17.52614769222031
$ws.Cells.Item(9, 7).Value = -17.52614769222031
$ws.Cells.Item(9, 8).Value = -17.52614769222031
$ws.Cells.Item(9, 9).Value = -17.52614769222031
$ws.Cells.Item(9, 10).Value = -17.52614769222031
$ws.Cells.Item(9, 11).Value = -17.52614769222031

# Row 10
$ws.Cells.Item(10, 2).Value = -17.52614769222031
$ws.Cells.Item(10, 3).Value = -17.52614769222031
$ws.Cells.Item(10, 4).Value = -17.52614769222031
$ws.Cells.Item(10, 5).Value = -17.52614769222031
$ws.Cells.Item(10, 6).Value = -17.52614769222031
$ws.Cells.Item(10, 7).Value = -17.52614769222031
$ws.Cells.Item(10, 8).Value = -17.52614769222031
$ws.Cells.Item(10, 9).Value = 0.1685162563634586
$ws.Cells.Item(10, 10).Value = -17.52614769222031
$ws.Cells.Item(10, 11).Value = 4.321920833776404

# Row 11
$ws.Cells.Item(11, 2).Value = -17.52614769222031
$ws.Cells.Item(11, 3).Value = -17.52614769222031
$ws.Cells.Item(11, 4).Value = -17.52614769222031
$ws.Cells.Item(11, 5).Value = 2.065310737360149
$ws.Cells.Item(11, 6).Value = -17.52614769222031
$ws.Cells.Item(11, 7).Value = 1.617058622130137
$ws.Cells.Item(11, 8).Value = -17.52614769222031
$ws.Cells.Item(11, 9).Value = -17.52614769222031
$ws.Cells.Item(11, 10).Value = -17.52614769222031
$ws.Cells.Item(11, 11).Value = -17.52614769222031

# Row 12
$ws.Cells.Item(12, 2).Value = -17.52614769222031
$ws.Cells.Item(12, 3).Value = -17.52614769222031
$ws.Cells.Item(12, 4).Value = -17.52614769222031
$ws.Cells.Item(12, 5).Value = -17.52614769222031
$ws.Cells.Item(12, 6).Value = -17.52614769222031
$ws.Cells.Item(12, 7).Value = -17.52614769222031
$ws.Cells.Item(12, 8).Value = -17.52614769222031
$ws.Cells.Item(12, 9).Value = -17.52614769222031
$ws.Cells.Item(12, 10).Value = -17.52614769222031
$ws.Cells.Item(12, 11).Value = -17.52614769222031

# Row 13
$ws.Cells.Item(13, 2).Value = -17.52614769222031
$ws.Cells.Item(13, 3).Value = -17.52614769222031
$ws.Cells.Item(13, 4).Value = -17.52614769222031
$ws.Cells.Item(13, 5).Value = 1.934854055039418
$ws.Cells.Item(13, 6).Value = -17.52614769222031
$ws.Cells.Item(13, 7).Value = -17.52614769222031
$ws.Cells.Item(13, 8).Value = -17.52614769222031
$ws.Cells.Item(13, 9).Value = -17.52614769222031
$ws.Cells.Item(13, 10).Value = 0.5779647167867314
$ws.Cells.Item(13, 11).Value = -17.52614769222031

# Row 14
$ws.Cells.Item(14, 2).Value = -17.52614769222031
$ws.Cells.Item(14, 3).Value = -17.52614769222031
$ws.Cells.Item(14, 4).Value = 0.8587582892675306
$ws.Cells.Item(14, 5).Value = -17.52614769222031
$ws.Cells.Item(14, 6).Value = -17.52614769222031
$ws.Cells.Item(14, 7).Value = -17.52614769222031
$ws.Cells.Item(14, 8).Value = -17.52614769222031
$ws.Cells.Item(14, 9).Value = -17.52614769222031
$ws.Cells.Item(14, 10).Value = -17.52614769222031
$ws.Cells.Item(14, 11).Value = -17.52614769222031

# Row 15
$ws.Cells.Item(15, 2).Value = -17.52614769222031
$ws.Cells.Item(15, 3).Value = -17.52614769222031
$ws.Cells.Item(15, 4).Value = -0.3080669738120843
$ws.Cells.Item(15, 5).Value = -17.52614769222031
$ws.Cells.Item(15, 6).Value = -17.52614769222031
$ws.Cells.Item(15, 7).Value = -17.52614769222031
$ws.Cells.Item(15, 8).Value = -17.52614769222031
$ws.Cells.Item(15, 9).Value = -17.52614769222031
$ws.Cells.Item(15, 10).Value = -17.52614769222031
$ws.Cells.Item(15, 11).Value = -17.52614769222031

# Row 16
$ws.Cells.Item(16, 2).Value = -17.52614769222031
$ws.Cells.Item(16, 3).Value = -17.52614769222031
$ws.Cells.Item(16, 4).Value = -17.52614769222031
$ws.Cells.Item(16, 5).Value = -17.52614769222031
$ws.Cells.Item(16, 6).Value = -17.52614769222031
$ws.Cells.Item(16, 7).Value = -17.52614769222031
$ws.Cells.Item(16, 8).Value = -17.52614769222031
$ws.Cells.Item(16, 9).Value = -17.52614769222031
$ws.Cells.Item(16, 10).Value = 2.150685750352566
$ws.Cells.Item(16, 11).Value = -17.52614769222031

# Row 17
$ws.Cells.Item(17, 2).Value = -17.52614769222031
$ws.Cells.Item(17, 3).Value = 0.327898862825529
$ws.Cells.Item(17, 4).Value = -0.1957822210714502
$ws.Cells.Item(17, 5).Value = -17.52614769222031
$ws.Cells.Item(17, 6).Value = -17.52614769222031
$ws.Cells.Item(17, 7).Value = -17.52614769222031
$ws.Cells.Item(17, 8).Value = 1.526306263826329
$ws.Cells.Item(17, 9).Value = -0.4246951586964121
$ws.Cells.Item(17, 10).Value = 1.773509756564114
$ws.Cells.Item(17, 11).Value = -17.52614769222031

# Row 18
$ws.Cells.Item(18, 2).Value = -17.52614769222031
$ws.Cells.Item(18, 3).Value = -17.52614769222031
$ws.Cells.Item(18, 4).Value = -17.52614769222031
$ws.Cells.Item(18, 5).Value = -17.52614769222031
$ws.Cells.Item(18, 6).Value = -17.52614769222031
$ws.Cells.Item(18, 7).Value = -17.52614769222031
$ws.Cells.Item(18, 8).Value = 2.895992034378907
$ws.Cells.Item(18, 9).Value = -0.2911338494485357
$ws.Cells.Item(18, 10).Value = 2.098171216974708
$ws.Cells.Item(18, 11).Value = -17.52614769222031

# Row 19
$ws.Cells.Item(19, 2).Value = -17.52614769222031
$ws.Cells.Item(19, 3).Value = -17.52614769222031
$ws.Cells.Item(19, 4).Value = 2.967525040241901
$ws.Cells.Item(19, 5).Value = -17.52614769222031
$ws.Cells.Item(19, 6).Value = -17.52614769222031
$ws.Cells.Item(19, 7).Value = -17.52614769222031
$ws.Cells.Item(19, 8).Value = 1.520510068357404
$ws.Cells.Item(19, 9).Value = 1.082486151322128
$ws.Cells.Item(19, 10).Value = -17.52614769222031
$ws.Cells.Item(19, 11).Value = -17.52614769222031

# Row 20
$ws.Cells.Item(20, 2).Value = -17.52614769222031
$ws.Cells.Item(20, 3).Value = 3.215856577220632
$ws.Cells.Item(20, 4).Value = 2.927276684615945
$ws.Cells.Item(20, 5).Value = -17.52614769222031
$ws.Cells.Item(20, 6).Value = 2.213107168763577
$ws.Cells.Item(20, 7).Value = -17.52614769222031
$ws.Cells.Item(20, 8).Value = 1.013741885935098
$ws.Cells.Item(20, 9).Value = 3.771740852984725
$ws.Cells.Item(20, 10).Value = -17.52614769222031
$ws.Cells.Item(20, 11).Value = -17.52614769222031

# Row 21
$ws.Cells.Item(21, 2).Value = -17.52614769222031
$ws.Cells.Item(21, 3).Value = 2.712850568109811
$ws.Cells.Item(21, 4).Value = -17.52614769222031
$ws.Cells.Item(21, 5).Value = 3.269518649336442
$ws.Cells.Item(21, 6).Value = -17.52614769222031
$ws.Cells.Item(21, 7).Value = 2.662287822290175
$ws.Cells.Item(21, 8).Value = 1.064764588184815
$ws.Cells.Item(21, 9).Value = -17.52614769222031
$ws.Cells.Item(21, 10).Value = -17.52614769222031
$ws.Cells.Item(21, 11).Value = -17.52614769222031

